# Update the Training Set Import Template header row.
# - D1 "response_impact_1" -> "slider_impact_direction"
# - E1 "response_impact_2" -> "choice_1_impact"
# - F1 "etc."              -> "choice_2_impact"
# - Add G1..N1 "choice_3_impact" .. "choice_10_impact" (matching header style)
# - Extend the custom column width (17.14) from D:E to also cover column F

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "slider_impact_direction"
$ws.Range("E1").Value = "choice_1_impact"
$ws.Range("F1").Value = "choice_2_impact"
$ws.Range("G1").Value = "choice_3_impact"
$ws.Range("H1").Value = "choice_4_impact"
$ws.Range("I1").Value = "choice_5_impact"
$ws.Range("J1").Value = "choice_6_impact"
$ws.Range("K1").Value = "choice_7_impact"
$ws.Range("L1").Value = "choice_8_impact"
$ws.Range("M1").Value = "choice_9_impact"
$ws.Range("N1").Value = "choice_10_impact"

# Give the newly added header cells the same formatting as the existing
# header cells (style used by A1:F1), by copying the format from F1.
$ws.Range("F1").Copy()
$ws.Range("G1:N1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Widen the custom D:E column width so it also covers the new column F.
$ws.Range("F1").EntireColumn.ColumnWidth = 16.33
